# Updated cryptos list on Fri Jun 16 13:24:09 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a leading apostrophe to force text entry (matches original inlineStr
# cells) so numeric-looking strings like "0.9984" are not coerced to
# numbers by Excel's auto-detection.

$ws.Range("D2").Value  = "'25.558.66"
$ws.Range("E2").Value  = "'  +2.36%  "

$ws.Range("D3").Value  = "'1.667.45"
$ws.Range("E3").Value  = "'  +1.78%  "

$ws.Range("D4").Value  = "'0.9984"
$ws.Range("E4").Value  = "'  +0.19%  "

$ws.Range("D5").Value  = "'235.19"
$ws.Range("E5").Value  = "'  +1.00%  "

$ws.Range("E6").Value  = "'  +0.05%  "

$ws.Range("D7").Value  = "'0.4650"
$ws.Range("E7").Value  = "'  -3.05%  "

$ws.Range("D8").Value  = "'0.2575"
$ws.Range("E8").Value  = "'  -0.89%  "

$ws.Range("D9").Value  = "'0.06129"
$ws.Range("E9").Value  = "'  +0.49%  "

$ws.Range("D10").Value = "'1.665.67"
$ws.Range("E10").Value = "'  +1.68%  "

$ws.Range("D11").Value = "'0.06949"
$ws.Range("E11").Value = "'  -1.56%  "

$ws.Range("E12").Value = "'  +0.66%  "

$ws.Range("D13").Value = "'4.344"
$ws.Range("E13").Value = "'  -0.70%  "

$ws.Range("D14").Value = "'74.89"
$ws.Range("E14").Value = "'  +1.96%  "

$ws.Range("D15").Value = "'0.5719"
$ws.Range("E15").Value = "'  -4.43%  "

$ws.Range("D16").Value = "'0.9997"
$ws.Range("E16").Value = "'  +0.08%  "

$ws.Range("D17").Value = "'0.9990"
$ws.Range("E17").Value = "'  +0.15%  "

$ws.Range("D18").Value = "'25.557.52"
$ws.Range("E18").Value = "'  +2.44%  "

$ws.Range("D19").Value = "'0.000006720"
$ws.Range("E19").Value = "'  +1.96%  "

$ws.Range("E20").Value = "'  +1.04%  "

$ws.Range("D21").Value = "'1.877.56"
$ws.Range("E21").Value = "'  +1.56%  "

$ws.Range("D22").Value = "'4.416"
$ws.Range("E22").Value = "'  +0.87%  "

$ws.Range("D23").Value = "'8.676"
$ws.Range("E23").Value = "'  +0.87%  "

$ws.Range("D24").Value = "'5.229"
$ws.Range("E24").Value = "'  -0.33%  "

$ws.Range("D25").Value = "'134.60"
$ws.Range("E25").Value = "'  +0.98%  "

$ws.Range("D26").Value = "'14.82"
$ws.Range("E26").Value = "'  -0.41%  "

$ws.Range("D27").Value = "'1.365"
$ws.Range("E27").Value = "'  -1.40%  "

$ws.Range("D28").Value = "'1.712"
$ws.Range("E28").Value = "'  +4.48%  "

$ws.Range("D29").Value = "'103.62"
$ws.Range("E29").Value = "'  -0.61%  "

$ws.Range("D30").Value = "'3.964"
$ws.Range("E30").Value = "'  +2.92%  "

$ws.Range("D31").Value = "'0.07699"
$ws.Range("E31").Value = "'  +0.12%  "

$ws.Range("D32").Value = "'3.599"
$ws.Range("E32").Value = "'  +1.58%  "

$ws.Range("D33").Value = "'0.04330"
$ws.Range("E33").Value = "'  +1.06%  "

$ws.Range("D34").Value = "'2.618"
$ws.Range("E34").Value = "'  +1.78%  "

$ws.Range("D35").Value = "'0.9437"
$ws.Range("E35").Value = "'  +1.80%  "

$ws.Range("D36").Value = "'0.5997"
$ws.Range("E36").Value = "'  +2.49%  "

$ws.Range("D37").Value = "'0.9137"
$ws.Range("E37").Value = "'  +10.10%  "

$ws.Range("D38").Value = "'2.483"
$ws.Range("E38").Value = "'  -2.53%  "

$ws.Range("D39").Value = "'105.47"
$ws.Range("E39").Value = "'  +7.12%  "

$ws.Range("D40").Value = "'0.9990"
$ws.Range("E40").Value = "'  +0.07%  "

# Rows 41 and 42 swap content: VeChain <-> RenderToken, each with new
# price/volume values.
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'1.834"
$ws.Range("E41").Value = "'  +5.11%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.01462"
$ws.Range("E42").Value = "'  -3.77%  "

$ws.Range("D43").Value = "'5.048"
$ws.Range("E43").Value = "'  +7.82%  "

$ws.Range("D44").Value = "'0.3706"
$ws.Range("E44").Value = "'  +0.31%  "

$ws.Range("E45").Value = "'  +2.14%  "

$ws.Range("D46").Value = "'0.05249"
$ws.Range("E46").Value = "'  +1.09%  "

$ws.Range("D47").Value = "'6.115"
$ws.Range("E47").Value = "'  +0.73%  "

$ws.Range("D48").Value = "'29.72"
$ws.Range("E48").Value = "'  +1.67%  "

$ws.Range("D49").Value = "'7.538"
$ws.Range("E49").Value = "'  +4.54%  "

$ws.Range("E50").Value = "'  +0.31%  "

$ws.Range("D51").Value = "'0.9992"
$ws.Range("E51").Value = "'  +0.31%  "
